$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 86, shifting existing rows 86-100 down to 87-101.
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row 86 with a duplicate of the (unchanged) row 85
# record, matching the target workbook state.
$ws.Range("A86").Value = 4
$ws.Range("B86").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C86").Value = "Los Lagos"
$ws.Range("D86").Value = 44726
$ws.Range("E86").Value = 10
$ws.Range("F86").Value = "Fruta"
$ws.Range("G86").Value = 100104
$ws.Range("H86").Value = "Frutos de pepita"
$ws.Range("I86").Value = 100104003
$ws.Range("J86").Value = "Membrillo"
$ws.Range("K86").Value = "Champion"
$ws.Range("L86").Value = "Primera"
$ws.Range("M86").Value = 300
$ws.Range("N86").Value = 13000
$ws.Range("O86").Value = 14000
$ws.Range("P86").Value = 13500
$ws.Range("Q86").Value = "`$/caja 18 kilos granel"
$ws.Range("R86").Value = "Región de O'Higgins"
$ws.Range("S86").Value = 750
$ws.Range("T86").Value = 18
